$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.127.95'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.87%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.772.95'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.86%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.009'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.70%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '333.47'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.73%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.004'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.61%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3765'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.91%  '

$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3399'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -2.77%  '

$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '48.12'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -3.00%  '

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.08%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07374'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -3.51%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.003'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.18%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.35'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.88%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.361'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.58%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.776.62'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.52%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.980'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -3.21%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001080'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -3.47%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06645'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.88%  '

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.95%  '

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.71%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.513'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.69%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '17.10'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -3.16%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.149.02'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.87%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.31'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -6.70%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.425'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.94%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.488'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.37%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.491'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.36%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.93'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +3.26%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '151.21'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.50%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.978.52'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.33%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '132.35'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.46%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.059'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.56%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.909'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -7.75%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08609'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.69%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '12.86'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -4.38%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.650'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -4.23%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.357'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -5.00%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.6750'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.47%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06279'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -3.47%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.02319'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -4.25%  '

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -4.38%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.663'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.65%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.230'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.94%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.32'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.43%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.004'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.64%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6281'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.80%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.826'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -3.57%  '

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.74%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '128.40'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.49%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07128'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -3.24%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '78.32'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.52%  '
